$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D7 (Water Use Efficiency) observation text
$ws.Range("D7").Value = "High and low producing nations have same water use efficiency. Need to further analyze which crops are the most water-efficient."

# Update the D8 (Agriculture share of Government Expenditure) observation text
$ws.Range("D8").Value = "Even with lower agriculture share of government expenditure, high-producing nations have more production per capita. This does not imply government should spend less on their agriculture sector."

# Row 8's wrapped text now needs two lines, so its height grows to fit
$ws.Rows.Item(8).RowHeight = 43.5

# Update the view: zoom and selection
$win = $excel.ActiveWindow
$win.Zoom = 95
$ws.Range("D7").Select()
